# "Generate Report for Handback"
#
# The handback-status report stamps the datetime every time the
# handback/handoff xliff report is (re)generated. Refresh the
# "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" cells that were produced by this run.
#
# These are plain text timestamps (stored as shared strings, formatted
# with a yyyy-mm-dd HH:mm:ss display style) rather than real date
# serials, so they are written back as literal strings.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for the first file row
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-24 11:07:27"

# --- zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-24 11:07:21"
$zhcn.Range("K2").Value = "2016-08-24 11:07:38"

# --- de-de sheet: Correspond Handback DateTime
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K2").Value = "2016-08-24 11:07:46"
